$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Insert a new column before column C on the "科系_櫃台_分機" sheet.
# The previous column C (extension numbers) shifts right to become column D.
$ws2.Columns("C:C").Insert()

# The new column C header duplicates the column D header text ("分機").
$ws2.Cells.Item(1, 3).Value = $ws2.Cells.Item(1, 4).Text

# Populate the new column C (rows 2-21) with a placeholder network
# exception value.
for ($r = 2; $r -le 21; $r++) {
    $ws2.Cells.Item($r, 3).Value = 22222
}

# Column B (rows 2-21) gets a Text number format applied.
$ws2.Range("B2:B21").NumberFormat = "@"

# Restore the selection/cursor position on each sheet.
$ws1.Range("B12").Select()
$ws2.Range("B7").Select()
